# Apply updated crypto price/volume values (GitHub Actions refresh).
# D-column cells that are plain decimal numbers must be force-written as
# text (leading apostrophe / quote-prefix) so Excel doesn't silently
# convert them to floating point numbers and mangle formats like
# trailing zeros (e.g. "0.780" -> 0.78).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.052.74"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.678.30"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'215.75"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  -1.14%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("D9").Value = "'21.37"
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("D10").Value = "'0.0623"
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "1.914.63"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "1.682.02"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'0.536"
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "'66.39"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "27.055.10"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").Value = "'8.16"
$ws.Range("E18").Value = "  +2.20%  "
$ws.Range("D19").Value = "'235.57"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "'4.46"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").Value = "'9.26"
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "'147.31"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'7.28"
$ws.Range("E26").Value = "  +1.87%  "
$ws.Range("D27").Value = "'16.52"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("D28").Value = "'0.113"
$ws.Range("E28").Value = "  -1.01%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'0.0497"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'3.38"
$ws.Range("E32").Value = "  +0.80%  "
$ws.Range("D33").Value = "1.543.42"
$ws.Range("E33").Value = "  +6.34%  "
$ws.Range("D34").Value = "'3.17"
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +4.31%  "
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'0.587"
$ws.Range("E37").Value = "  +0.76%  "
$ws.Range("D38").Value = "'0.913"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("D40").Value = "'1.05"
$ws.Range("E40").Value = "  +6.58%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("D42").Value = "'67.87"
$ws.Range("D43").Value = "'5.54"
$ws.Range("E43").Value = "  -3.62%  "
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "1.822.05"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'0.780"
$ws.Range("E46").Value = "  -0.45%  "
$ws.Range("D47").Value = "'90.31"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "0.0₆0108"
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +1.38%  "
$ws.Range("D51").Value = "'8.02"
$ws.Range("E51").Value = "  +6.20%  "
